# ------------------------------------------------------------------
# 1. Refresh the "time_taken" timestamps on the existing "data" sheet
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

$timestamps = @(
    "2021-10-05 14:19:47.306321",
    "2021-10-05 14:19:47.306329",
    "2021-10-05 14:19:47.306332",
    "2021-10-05 14:19:47.306335",
    "2021-10-05 14:19:47.306338",
    "2021-10-05 14:19:47.306340",
    "2021-10-05 14:19:47.306343",
    "2021-10-05 14:19:47.306345",
    "2021-10-05 14:19:47.306348",
    "2021-10-05 14:19:47.306350",
    "2021-10-05 14:19:47.306353",
    "2021-10-05 14:19:47.306355",
    "2021-10-05 14:19:47.306358",
    "2021-10-05 14:19:47.306360",
    "2021-10-05 14:19:47.306363",
    "2021-10-05 14:19:47.306365",
    "2021-10-05 14:19:47.306368",
    "2021-10-05 14:19:47.306370",
    "2021-10-05 14:19:47.306373",
    "2021-10-05 14:19:47.306375",
    "2021-10-05 14:19:47.306378",
    "2021-10-05 14:19:47.306380",
    "2021-10-05 14:19:47.306383"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataWs.Cells.Item($row, 6).Value = $timestamps[$i]
}

# ------------------------------------------------------------------
# 2. Add the new "metadata" sheet right after "data"
# ------------------------------------------------------------------
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

# Header row
$metaWs.Cells.Item(1, 2).Value = "data_name"
$metaWs.Cells.Item(1, 3).Value = "data_id"
$metaWs.Cells.Item(1, 4).Value = "data_version"
$metaWs.Cells.Item(1, 5).Value = "data_version_created"
$metaWs.Cells.Item(1, 6).Value = "panel_query_time"
$metaWs.Cells.Item(1, 7).Value = "panel_get_request"

$headerRng = $metaWs.Range("B1:G1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# Data row
$metaWs.Cells.Item(2, 1).Value = 0
$a2 = $metaWs.Cells.Item(2, 1)
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$metaWs.Cells.Item(2, 2).Value = "Corneal dystrophies"
$metaWs.Cells.Item(2, 3).Value = 658
$metaWs.Cells.Item(2, 4).NumberFormat = "@"
$metaWs.Cells.Item(2, 4).Value = "1.6"
$metaWs.Cells.Item(2, 4).NumberFormat = "General"
$metaWs.Cells.Item(2, 5).Value = "2020-12-15T10:59:24.518016Z"
$metaWs.Cells.Item(2, 6).Value = "2021-10-05 14:19:47.302594"
$metaWs.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/658/?format=json"

$dataWs.Activate()
